# Update BOM for SPEEDY

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MCU part number revision (new shared string created first)
$ws.Range("A18").Value = "STM32F303CCT6"

# Resistor/value comments simplified from "<x>Ohm" style to plain "<x>k"/number style
$ws.Range("A6").Value = "10k"
$ws.Range("A8").Value = "4.7k"
$ws.Range("A11").Value = "1.5k"
$ws.Range("A13").Value = "100k"
$ws.Range("A15").Value = "1k"

# Plain-number comment cells (no longer stored as text)
$ws.Range("A5").Value = 22
$ws.Range("A12").Value = 100

# Clear the now-unused "LCSC Part #(optional)" values in column D
$ws.Range("D2").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("D5").Value = ""

# Removed rows worth of parts (switch / connector) - clear their comment/designator
$ws.Range("A14").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("A16").Value = ""
$ws.Range("B16").Value = ""

# Voltage regulator part number revision (new shared string created last)
$ws.Range("A7").Value = "me6211-3.3V"

# Update the active selection to match the saved file
$ws.Range("A8").Select()
